$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Junio de 2020 a las 17:47"

# Row 4
$ws.Cells.Item(4, 2).Value = 2339196
$ws.Cells.Item(4, 3).Value = 8618
$ws.Cells.Item(4, 4).Value = 974288
$ws.Cells.Item(4, 5).Value = 1242863
$ws.Cells.Item(4, 7).Value = 65
$ws.Cells.Item(4, 8).Value = 122045

# Row 7
$ws.Cells.Item(7, 2).Value = 421279
$ws.Cells.Item(7, 3).Value = 9552
$ws.Cells.Item(7, 4).Value = 234714
$ws.Cells.Item(7, 5).Value = 173077
$ws.Cells.Item(7, 7).Value = 211
$ws.Cells.Item(7, 8).Value = 13488

# Row 14
$ws.Cells.Item(14, 2).Value = 191276
$ws.Cells.Item(14, 3).Value = 60
$ws.Cells.Item(14, 5).Value = 7415

# Row 21
$ws.Cells.Item(21, 2).Value = 101286
$ws.Cells.Item(21, 3).Value = 267
$ws.Cells.Item(21, 4).Value = 63860
$ws.Cells.Item(21, 5).Value = 28996
$ws.Cells.Item(21, 7).Value = 20
$ws.Cells.Item(21, 8).Value = 8430

# Row 45
$ws.Cells.Item(45, 2).Value = 26677
$ws.Cells.Item(45, 3).Value = 899
$ws.Cells.Item(45, 4).Value = 15138
$ws.Cells.Item(45, 5).Value = 10877
$ws.Cells.Item(45, 7).Value = 7
$ws.Cells.Item(45, 8).Value = 662

# Row 51
$ws.Cells.Item(51, 2).Value = 20734
$ws.Cells.Item(51, 3).Value = 101
$ws.Cells.Item(51, 4).Value = 15685
$ws.Cells.Item(51, 5).Value = 4743
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 306

# Row 56
$ws.Cells.Item(56, 4).Value = 10897
$ws.Cells.Item(56, 5).Value = 6210

# Row 57
$ws.Cells.Item(57, 2).Value = 14200
$ws.Cells.Item(57, 3).Value = 247
$ws.Cells.Item(57, 5).Value = 5831
$ws.Cells.Item(57, 7).Value = 9
$ws.Cells.Item(57, 8).Value = 473

# Row 67
$ws.Cells.Item(67, 2).Value = 10463
$ws.Cells.Item(67, 3).Value = 15
$ws.Cells.Item(67, 4).Value = 7498
$ws.Cells.Item(67, 5).Value = 2628
$ws.Cells.Item(67, 7).Value = 1
$ws.Cells.Item(67, 8).Value = 337

# Row 84
$ws.Cells.Item(84, 2).Value = 4582
$ws.Cells.Item(84, 3).Value = 17
$ws.Cells.Item(84, 4).Value = 3859
$ws.Cells.Item(84, 5).Value = 678

# Row 88
$ws.Cells.Item(88, 2).Value = 4120
$ws.Cells.Item(88, 3).Value = 15
$ws.Cells.Item(88, 4).Value = 3956
$ws.Cells.Item(88, 5).Value = 54

# Row 93
$ws.Cells.Item(93, 2).Value = 3266
$ws.Cells.Item(93, 3).Value = 10
$ws.Cells.Item(93, 5).Value = 1702

# Row 127
$ws.Cells.Item(127, 2).Value = 1033
$ws.Cells.Item(127, 3).Value = 18
$ws.Cells.Item(127, 4).Value = 739
$ws.Cells.Item(127, 5).Value = 285

# Row 140
$ws.Cells.Item(140, 1).Value = "Mozambique"
$ws.Cells.Item(140, 2).Value = 733
$ws.Cells.Item(140, 3).Value = 45
$ws.Cells.Item(140, 4).Value = 181
$ws.Cells.Item(140, 5).Value = 547
$ws.Cells.Item(140, 7).Value = 1
$ws.Cells.Item(140, 8).Value = 5

# Row 141
$ws.Cells.Item(141, 1).Value = "Malaui"
$ws.Cells.Item(141, 2).Value = 730
$ws.Cells.Item(141, 3).Value = 110
$ws.Cells.Item(141, 4).Value = 258
$ws.Cells.Item(141, 5).Value = 461
$ws.Cells.Item(141, 7).Value = 3
$ws.Cells.Item(141, 8).Value = 11

# Row 142
$ws.Cells.Item(142, 1).Value = "Crucero"
$ws.Cells.Item(142, 2).Value = 712
$ws.Cells.Item(142, 4).Value = 651
$ws.Cells.Item(142, 5).Value = 48
$ws.Cells.Item(142, 8).Value = 13

# Row 143
$ws.Cells.Item(143, 1).Value = "Ruanda"
$ws.Cells.Item(143, 2).Value = 702
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 357
$ws.Cells.Item(143, 5).Value = 343
$ws.Cells.Item(143, 8).Value = 2

# Row 144
$ws.Cells.Item(144, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(144, 2).Value = 698
$ws.Cells.Item(144, 3).Value = 5
$ws.Cells.Item(144, 4).Value = 203
$ws.Cells.Item(144, 5).Value = 483
$ws.Cells.Item(144, 8).Value = 12

# Row 145
$ws.Cells.Item(145, 1).Value = "San Marino"
$ws.Cells.Item(145, 2).Value = 696
$ws.Cells.Item(145, 4).Value = 610
$ws.Cells.Item(145, 5).Value = 44
$ws.Cells.Item(145, 8).Value = 42

# Row 149
$ws.Cells.Item(149, 2).Value = 626
$ws.Cells.Item(149, 3).Value = 25
$ws.Cells.Item(149, 5).Value = 338
$ws.Cells.Item(149, 7).Value = 1
$ws.Cells.Item(149, 8).Value = 34

# Row 202
$ws.Cells.Item(202, 1).Value = "Fiyi"

# Row 203
$ws.Cells.Item(203, 1).Value = "Dominica"

# Row 208
$ws.Cells.Item(208, 1).Value = "Santa Sede"
$ws.Cells.Item(208, 4).Value = 12
$ws.Cells.Item(208, 8).Value = 0

# Row 209
$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 8).Value = 1

# Row 213
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 8).Value = 0

# Row 214
$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 8).Value = 1
